$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: establish styles for the brand-new rows 279-281 by copying
#            format from existing template rows (id col A = bold/border style,
#            date col E = datetime style). Done first, before source rows'
#            values are overwritten below. ---
$ws.Range("A278").Copy($ws.Range("A279"))
$ws.Range("E278").Copy($ws.Range("E279"))
$ws.Range("A277").Copy($ws.Range("A280"))
$ws.Range("E277").Copy($ws.Range("E280"))
$ws.Range("A278").Copy($ws.Range("A281"))
$ws.Range("E278").Copy($ws.Range("E281"))

# --- Step 2: rows 235-239 -- cyclic shuffle of the B:AC betting-odds data
#            (columns A id and E date are unchanged) ---
# Row 235
$ws.Range("B235").Value = 6861095
$ws.Range("C235").Value = "Romania Liga I"
$ws.Range("D235").Value = "Romania Liga I"
$ws.Range("F235").Value = "FC Botosani"
$ws.Range("G235").Value = "Farul Constanta"
$ws.Range("H235").Value = 0
$ws.Range("I235").Value = 0
$ws.Range("J235").Value = "D"
$ws.Range("K235").Value = 3.75
$ws.Range("L235").Value = 3.4
$ws.Range("M235").Value = 1.909
$ws.Range("N235").Value = 3.1
$ws.Range("O235").Value = 3
$ws.Range("P235").Value = 2.375
$ws.Range("Q235").Value = 0.25
$ws.Range("R235").Value = 1.775
$ws.Range("S235").Value = 2.1
$ws.Range("T235").Value = 2
$ws.Range("U235").Value = 1.8
$ws.Range("V235").Value = 2.05
$ws.Range("W235").Value = -1
$ws.Range("X235").Value = 2
$ws.Range("Y235").Value = -1
$ws.Range("Z235").Value = 0.3875
$ws.Range("AA235").Value = -0.5
$ws.Range("AB235").Value = -1
$ws.Range("AC235").Value = 1.05

# Row 236
$ws.Range("B236").Value = 6870268
$ws.Range("C236").Value = "Romania Liga I"
$ws.Range("D236").Value = "Romania Liga I"
$ws.Range("F236").Value = "Petrolul Ploiesti"
$ws.Range("G236").Value = "ACS Sepsi"
$ws.Range("H236").Value = 1
$ws.Range("I236").Value = 2
$ws.Range("J236").Value = "A"
$ws.Range("K236").Value = 2.8
$ws.Range("L236").Value = 3
$ws.Range("M236").Value = 2.55
$ws.Range("N236").Value = 3
$ws.Range("O236").Value = 3.2
$ws.Range("P236").Value = 2.3
$ws.Range("Q236").Value = 0.25
$ws.Range("R236").Value = 1.85
$ws.Range("S236").Value = 2
$ws.Range("T236").Value = 2.25
$ws.Range("U236").Value = 1.875
$ws.Range("V236").Value = 1.975
$ws.Range("W236").Value = -1
$ws.Range("X236").Value = -1
$ws.Range("Y236").Value = 1.3
$ws.Range("Z236").Value = -1
$ws.Range("AA236").Value = 1
$ws.Range("AB236").Value = 0.875
$ws.Range("AC236").Value = -1

# Row 237
$ws.Range("B237").Value = 6865915
$ws.Range("C237").Value = "Romania Liga I"
$ws.Range("D237").Value = "Romania Liga I"
$ws.Range("F237").Value = "FC Voluntari"
$ws.Range("G237").Value = "Universitatea Cluj"
$ws.Range("H237").Value = 0
$ws.Range("I237").Value = 0
$ws.Range("J237").Value = "D"
$ws.Range("K237").Value = 3.5
$ws.Range("L237").Value = 3.25
$ws.Range("M237").Value = 2.05
$ws.Range("N237").Value = 3.4
$ws.Range("O237").Value = 3.1
$ws.Range("P237").Value = 2.15
$ws.Range("Q237").Value = 0.25
$ws.Range("R237").Value = 1.975
$ws.Range("S237").Value = 1.875
$ws.Range("T237").Value = 2.25
$ws.Range("U237").Value = 2.05
$ws.Range("V237").Value = 1.75
$ws.Range("W237").Value = -1
$ws.Range("X237").Value = 2.1
$ws.Range("Y237").Value = -1
$ws.Range("Z237").Value = 0.4875
$ws.Range("AA237").Value = -0.5
$ws.Range("AB237").Value = -1
$ws.Range("AC237").Value = 0.75

# Row 238
$ws.Range("B238").Value = 6836277
$ws.Range("C238").Value = "Romania Liga I"
$ws.Range("D238").Value = "Romania Liga I"
$ws.Range("F238").Value = "CFR Cluj"
$ws.Range("G238").Value = "AFC Hermannstadt"
$ws.Range("H238").Value = 1
$ws.Range("I238").Value = 0
$ws.Range("J238").Value = "H"
$ws.Range("K238").Value = 1.7
$ws.Range("L238").Value = 3.4
$ws.Range("M238").Value = 5
$ws.Range("N238").Value = 1.65
$ws.Range("O238").Value = 3.5
$ws.Range("P238").Value = 5.25
$ws.Range("Q238").Value = -0.75
$ws.Range("R238").Value = 1.85
$ws.Range("S238").Value = 2
$ws.Range("T238").Value = 2.25
$ws.Range("U238").Value = 1.875
$ws.Range("V238").Value = 1.975
$ws.Range("W238").Value = 0.6499999999999999
$ws.Range("X238").Value = -1
$ws.Range("Y238").Value = -1
$ws.Range("Z238").Value = 0.425
$ws.Range("AA238").Value = -0.5
$ws.Range("AB238").Value = -1
$ws.Range("AC238").Value = 0.9750000000000001

# Row 239
$ws.Range("B239").Value = 6852370
$ws.Range("C239").Value = "Romania Liga I"
$ws.Range("D239").Value = "Romania Liga I"
$ws.Range("F239").Value = "Dinamo Bucharest"
$ws.Range("G239").Value = "ACS UTA Batrana Doamna"
$ws.Range("H239").Value = 1
$ws.Range("I239").Value = 0
$ws.Range("J239").Value = "H"
$ws.Range("K239").Value = 2.55
$ws.Range("L239").Value = 2.875
$ws.Range("M239").Value = 3
$ws.Range("N239").Value = 2.375
$ws.Range("O239").Value = 3
$ws.Range("P239").Value = 3.1
$ws.Range("Q239").Value = -0.25
$ws.Range("R239").Value = 2
$ws.Range("S239").Value = 1.85
$ws.Range("T239").Value = 2.25
$ws.Range("U239").Value = 1.975
$ws.Range("V239").Value = 1.875
$ws.Range("W239").Value = 1.375
$ws.Range("X239").Value = -1
$ws.Range("Y239").Value = -1
$ws.Range("Z239").Value = 1
$ws.Range("AA239").Value = -1
$ws.Range("AB239").Value = -1
$ws.Range("AC239").Value = 0.875

# --- Step 3: rows 277-278 -- replaced with newly-scraped match data
#            (adds H/I/J score/result + AB/AC columns that were blank before) ---
# Row 277
$ws.Range("A277").Value = 275
$ws.Range("B277").Value = 7951796
$ws.Range("C277").Value = "Romania Liga I"
$ws.Range("D277").Value = "Romania Liga I"
$ws.Range("E277").Value = 45402.375
$ws.Range("F277").Value = "CSM Politehnica Iasi"
$ws.Range("G277").Value = "FC Voluntari"
$ws.Range("H277").Value = 3
$ws.Range("I277").Value = 1
$ws.Range("J277").Value = "H"
$ws.Range("K277").Value = 2.25
$ws.Range("L277").Value = 3
$ws.Range("M277").Value = 3.25
$ws.Range("N277").Value = 2.45
$ws.Range("O277").Value = 3
$ws.Range("P277").Value = 2.875
$ws.Range("Q277").Value = 0
$ws.Range("R277").Value = 1.775
$ws.Range("S277").Value = 2.1
$ws.Range("T277").Value = 2.25
$ws.Range("U277").Value = 2.025
$ws.Range("V277").Value = 1.825
$ws.Range("W277").Value = 1.45
$ws.Range("X277").Value = -1
$ws.Range("Y277").Value = -1
$ws.Range("Z277").Value = 0.7749999999999999
$ws.Range("AA277").Value = -1
$ws.Range("AB277").Value = 1.025
$ws.Range("AC277").Value = -1

# Row 278
$ws.Range("A278").Value = 276
$ws.Range("B278").Value = 7951797
$ws.Range("C278").Value = "Romania Liga I"
$ws.Range("D278").Value = "Romania Liga I"
$ws.Range("E278").Value = 45402.47916666666
$ws.Range("F278").Value = "FC U Craiova 1948"
$ws.Range("G278").Value = "Dinamo Bucharest"
$ws.Range("H278").Value = 1
$ws.Range("I278").Value = 1
$ws.Range("J278").Value = "D"
$ws.Range("K278").Value = 2.2
$ws.Range("L278").Value = 3.1
$ws.Range("M278").Value = 3.25
$ws.Range("N278").Value = 2.45
$ws.Range("O278").Value = 3
$ws.Range("P278").Value = 2.875
$ws.Range("Q278").Value = 0
$ws.Range("R278").Value = 1.775
$ws.Range("S278").Value = 2.1
$ws.Range("T278").Value = 2.25
$ws.Range("U278").Value = 2.05
$ws.Range("V278").Value = 1.8
$ws.Range("W278").Value = -1
$ws.Range("X278").Value = 2
$ws.Range("Y278").Value = -1
$ws.Range("Z278").Value = 0
$ws.Range("AA278").Value = -0
$ws.Range("AB278").Value = -0.5
$ws.Range("AC278").Value = 0.4

# --- Step 4: row 279 -- brand new match row ---
# Row 279
$ws.Range("A279").Value = 277
$ws.Range("B279").Value = 7951757
$ws.Range("C279").Value = "Romania Liga I"
$ws.Range("D279").Value = "Romania Liga I"
$ws.Range("E279").Value = 45402.60416666666
$ws.Range("F279").Value = "FCSB"
$ws.Range("G279").Value = "Rapid Bucuresti"
$ws.Range("H279").Value = 2
$ws.Range("I279").Value = 2
$ws.Range("J279").Value = "D"
$ws.Range("K279").Value = 1.85
$ws.Range("L279").Value = 3.5
$ws.Range("M279").Value = 4
$ws.Range("N279").Value = 1.909
$ws.Range("O279").Value = 3.6
$ws.Range("P279").Value = 3.75
$ws.Range("Q279").Value = -0.5
$ws.Range("R279").Value = 1.9
$ws.Range("S279").Value = 1.95
$ws.Range("T279").Value = 2.5
$ws.Range("U279").Value = 1.85
$ws.Range("V279").Value = 2
$ws.Range("W279").Value = -1
$ws.Range("X279").Value = 2.6
$ws.Range("Y279").Value = -1
$ws.Range("Z279").Value = -1
$ws.Range("AA279").Value = 0.95
$ws.Range("AB279").Value = 0.8500000000000001
$ws.Range("AC279").Value = -1

# --- Step 5: rows 280-281 -- new rows holding the data that used to live
#            in rows 277-278 (not-yet-played fixtures, no score/result yet) ---
# Row 280
$ws.Range("A280").Value = 278
$ws.Range("B280").Value = 7951759
$ws.Range("C280").Value = "Romania Liga I"
$ws.Range("D280").Value = "Romania Liga I"
$ws.Range("E280").Value = 45403.54166666666
$ws.Range("F280").Value = "ACS Sepsi"
$ws.Range("G280").Value = "CS U Craiova"
$ws.Range("K280").Value = 2.6
$ws.Range("L280").Value = 3.2
$ws.Range("M280").Value = 2.6
$ws.Range("N280").Value = 2.45
$ws.Range("O280").Value = 3.2
$ws.Range("P280").Value = 2.75
$ws.Range("Q280").Value = 0
$ws.Range("R280").Value = 1.8
$ws.Range("S280").Value = 2.05
$ws.Range("T280").Value = 2.25
$ws.Range("U280").Value = 1.8
$ws.Range("V280").Value = 2.05
$ws.Range("W280").Value = 0
$ws.Range("X280").Value = 0
$ws.Range("Y280").Value = 0
$ws.Range("Z280").Value = 0
$ws.Range("AA280").Value = 0

# Row 281
$ws.Range("A281").Value = 279
$ws.Range("B281").Value = 7951794
$ws.Range("C281").Value = "Romania Liga I"
$ws.Range("D281").Value = "Romania Liga I"
$ws.Range("E281").Value = 45404.54166666666
$ws.Range("F281").Value = "Petrolul Ploiesti"
$ws.Range("G281").Value = "ACS UTA Batrana Doamna"
$ws.Range("K281").Value = 2.1
$ws.Range("L281").Value = 3.1
$ws.Range("M281").Value = 3.4
$ws.Range("N281").Value = 2.05
$ws.Range("O281").Value = 3.1
$ws.Range("P281").Value = 3.5
$ws.Range("Q281").Value = -0.25
$ws.Range("R281").Value = 1.8
$ws.Range("S281").Value = 2.05
$ws.Range("T281").Value = 2.25
$ws.Range("U281").Value = 2
$ws.Range("V281").Value = 1.85
$ws.Range("W281").Value = 0
$ws.Range("X281").Value = 0
$ws.Range("Y281").Value = 0
$ws.Range("Z281").Value = 0
$ws.Range("AA281").Value = 0
